{"js": "// Update the date line and every \"AxB=C\" multiplication-fact cell in the\n// table to the new values for this generated worksheet output.\nconst replacements = [\n  [\"2024-12-13 Friday\", \"2024-12-14 Saturday\"],\n  [\"17\u00d767=1139\", \"99\u00d738=3762\"],\n  [\"41\u00d755=2255\", \"64\u00d744=2816\"],\n  [\"29\u00d756=1624\", \"15\u00d746=690\"],\n  [\"58\u00d754=3132\", \"34\u00d782=2788\"],\n  [\"96\u00d737=3552\", \"92\u00d733=3036\"],\n  [\"34\u00d734=1156\", \"80\u00d722=1760\"],\n  [\"50\u00d739=1950\", \"61\u00d758=3538\"],\n  [\"75\u00d788=6600\", \"52\u00d778=4056\"],\n  [\"18\u00d770=1260\", \"42\u00d775=3150\"],\n  [\"47\u00d778=3666\", \"84\u00d768=5712\"],\n  [\"82\u00d796=7872\", \"89\u00d753=4717\"],\n  [\"57\u00d747=2679\", \"99\u00d786=8514\"],\n  [\"58\u00d726=1508\", \"50\u00d765=3250\"],\n  [\"87\u00d770=6090\", \"61\u00d745=2745\"],\n  [\"36\u00d758=2088\", \"45\u00d722=990\"],\n  [\"19\u00d762=1178\", \"66\u00d786=5676\"],\n  [\"70\u00d799=6930\", \"20\u00d780=1600\"],\n  [\"45\u00d735=1575\", \"36\u00d793=3348\"],\n  [\"99\u00d759=5841\", \"45\u00d724=1080\"],\n  [\"79\u00d725=1975\", \"73\u00d779=5767\"],\n  [\"23\u00d759=1357\", \"35\u00d759=2065\"],\n  [\"12\u00d791=1092\", \"68\u00d724=1632\"],\n  [\"58\u00d775=4350\", \"38\u00d796=3648\"],\n  [\"86\u00d777=6622\", \"21\u00d745=945\"],\n  [\"97\u00d792=8924\", \"98\u00d748=4704\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n", "ps1": "# Update the worksheet date and every \"AxB=C\" multiplication-fact cell\n# in the table to the new values for this generated output.\n$d = $word.ActiveDocument\n\n$oldTexts = @(\"2024-12-13 Friday\", \"17\u00d767=1139\", \"41\u00d755=2255\", \"29\u00d756=1624\", \"58\u00d754=3132\", \"96\u00d737=3552\", \"34\u00d734=1156\", \"50\u00d739=1950\", \"75\u00d788=6600\", \"18\u00d770=1260\", \"47\u00d778=3666\", \"82\u00d796=7872\", \"57\u00d747=2679\", \"58\u00d726=1508\", \"87\u00d770=6090\", \"36\u00d758=2088\", \"19\u00d762=1178\", \"70\u00d799=6930\", \"45\u00d735=1575\", \"99\u00d759=5841\", \"79\u00d725=1975\", \"23\u00d759=1357\", \"12\u00d791=1092\", \"58\u00d775=4350\", \"86\u00d777=6622\", \"97\u00d792=8924\")\n$newTexts = @(\"2024-12-14 Saturday\", \"99\u00d738=3762\", \"64\u00d744=2816\", \"15\u00d746=690\", \"34\u00d782=2788\", \"92\u00d733=3036\", \"80\u00d722=1760\", \"61\u00d758=3538\", \"52\u00d778=4056\", \"42\u00d775=3150\", \"84\u00d768=5712\", \"89\u00d753=4717\", \"99\u00d786=8514\", \"50\u00d765=3250\", \"61\u00d745=2745\", \"45\u00d722=990\", \"66\u00d786=5676\", \"20\u00d780=1600\", \"36\u00d793=3348\", \"45\u00d724=1080\", \"73\u00d779=5767\", \"35\u00d759=2065\", \"68\u00d724=1632\", \"38\u00d796=3648\", \"21\u00d745=945\", \"98\u00d748=4704\")\n\nfor ($i = 0; $i -lt $oldTexts.Count; $i++) {\n    $oldText = $oldTexts[$i]\n    $newText = $newTexts[$i]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
